# Update investment-capacity results with fresh values from server run.
$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0.1850090999999999
$ws.Range("B2").Value = 0.08791902848523356
$ws.Range("E2").Value = 0.2072241724204126
$ws.Range("G2").Value = 0.1304011109248467
$ws.Range("I2").Value = 0.9100542383090227
$ws.Range("M2").Value = 0.08588256169097744
$ws.Range("N2").Value = 8.387250548205344
$ws.Range("O2").Value = 5.572860925603398

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("E2").Value = 0.2277642275795874
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = 0.4837123116909772
$ws.Range("M2").Value = 0.04498818830902257
$ws.Range("N2").Value = 9.055699121367166
$ws.Range("O2").Value = 2.957574604300432

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("E2").Value = 0.2720621710968099
$ws.Range("G2").Value = 0.119088125212652
$ws.Range("I2").Value = 0.4493937760357649
$ws.Range("L2").Value = 0.1230520827199775
$ws.Range("M2").Value = 0.05330024124425767
$ws.Range("N2").Value = 5.054814646671876
$ws.Range("O2").Value = 2.372323717542278
